# Apply the edits described by the commit "Finished parameter estimates,
# started on quantiles" to the active worksheet.
#
# The changes are:
#   - Updated several parameter-estimate inputs in the "Summer" / "Fall"
#     blocks (columns B:E) that feed the BIC/AIC computations in columns
#     G:K.
#   - Filled in two previously-empty columns (E7:E8 and E19:E20) that are
#     referenced by the AIC block formulas (which used to evaluate to 0).
#   - Removed the stray "**" footnote marker that lived in D10/J10 (and,
#     as a consequence of Excel's shared-string cleanup, the now-unused
#     "**" shared string disappears from the workbook - this happens
#     automatically, we just need to clear the cells).
#   - Moved the active cell selection from F8 to F11.
#
# All of the other cells touched by the diff are plain formulas
# (=B3*2+B4*LN(24*24*2), =C3, etc.) that already exist in the workbook and
# simply recalculate once the inputs below change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: Summer, rows 3-4 (neg log-lik / df) ---------------------
$ws.Range("C3").Value = 558.97
$ws.Range("D3").Value = 445.39

$ws.Range("C4").Value = 174.37
$ws.Range("D4").Value = 215.91

# --- Block 2: Fall, rows 7-8 (neg log-lik / df) ------------------------
$ws.Range("D7").Value = 666.18
$ws.Range("E7").Value = 499.69

$ws.Range("D8").Value = 167.24
$ws.Range("E8").Value = 226.4

# --- Remove the "**" footnote markers in row 10 ------------------------
$ws.Range("D10").ClearContents()
$ws.Range("J10").ClearContents()

# --- Block 3: Summer (bl+dfbl), row 16 (df) ----------------------------
$ws.Range("C16").Value = 57.22

# --- Block 4: Fall (bl+dfbl), rows 19-20 (neg log-lik / df) ------------
$ws.Range("D19").Value = -1552.77
$ws.Range("E19").Value = -1792.03

$ws.Range("D20").Value = 20.88
$ws.Range("E20").Value = 23.09

# --- Update the saved active-cell selection ----------------------------
$ws.Range("F11").Select()

$wb.Save()
